$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing DescEmb-RNN_Scr -> DescEmb-RNN_Transfer tags (rows 454-456, 497) ---
$ws.Range("Q454").Value = "DescEmb-RNN_Transfer"
$ws.Range("Q455").Value = "DescEmb-RNN_Transfer"
$ws.Range("Q456").Value = "DescEmb-RNN_Transfer"
$ws.Range("Q497").Value = "DescEmb-RNN_Transfer"

# --- Append new experiment rows 499-507 (DSVA_DPE DescEmb-RNN transfer-learning runs) ---
$ws.Cells.Item(499, 1).Value = "outputs/2024-05-06/05-12-10"
$ws.Cells.Item(499, 2).Value = $true
$ws.Cells.Item(499, 3).Value = "mimiciii"
$ws.Cells.Item(499, 4).Value = "los_3day"
$ws.Cells.Item(499, 5).Value = "DSVA_DPE"
$ws.Cells.Item(499, 6).Value = "descemb_rnn"
$ws.Cells.Item(499, 7).Value = "ehr_model"
$ws.Cells.Item(499, 8).Value = $false
$ws.Cells.Item(499, 9).Value = $false
$ws.Cells.Item(499, 10).Value = $false
$ws.Cells.Item(499, 11).Value = $true
$ws.Cells.Item(499, 13).Value = 75
$ws.Cells.Item(499, 14).Value = 1.519
$ws.Cells.Item(499, 15).Value = 0.51
$ws.Cells.Item(499, 16).Value = 0.349
$ws.Cells.Item(499, 17).Value = "DescEmb-RNN_Transfer"

$ws.Cells.Item(500, 1).Value = "outputs/2024-05-06/05-39-05"
$ws.Cells.Item(500, 2).Value = $true
$ws.Cells.Item(500, 3).Value = "mimiciii"
$ws.Cells.Item(500, 4).Value = "los_7day"
$ws.Cells.Item(500, 5).Value = "DSVA_DPE"
$ws.Cells.Item(500, 6).Value = "descemb_rnn"
$ws.Cells.Item(500, 7).Value = "ehr_model"
$ws.Cells.Item(500, 8).Value = $false
$ws.Cells.Item(500, 9).Value = $false
$ws.Cells.Item(500, 10).Value = $false
$ws.Cells.Item(500, 11).Value = $true
$ws.Cells.Item(500, 13).Value = 80
$ws.Cells.Item(500, 14).Value = 1.178
$ws.Cells.Item(500, 15).Value = 0.512
$ws.Cells.Item(500, 16).Value = 0.134
$ws.Cells.Item(500, 17).Value = "DescEmb-RNN_Transfer"

$ws.Cells.Item(501, 1).Value = "outputs/2024-05-06/06-07-43"
$ws.Cells.Item(501, 2).Value = $true
$ws.Cells.Item(501, 3).Value = "mimiciii"
$ws.Cells.Item(501, 4).Value = "readmission"
$ws.Cells.Item(501, 5).Value = "DSVA_DPE"
$ws.Cells.Item(501, 6).Value = "descemb_rnn"
$ws.Cells.Item(501, 7).Value = "ehr_model"
$ws.Cells.Item(501, 8).Value = $false
$ws.Cells.Item(501, 9).Value = $false
$ws.Cells.Item(501, 10).Value = $false
$ws.Cells.Item(501, 11).Value = $true
$ws.Cells.Item(501, 13).Value = 113
$ws.Cells.Item(501, 14).Value = 0.855
$ws.Cells.Item(501, 15).Value = 0.503
$ws.Cells.Item(501, 16).Value = 0.043
$ws.Cells.Item(501, 17).Value = "DescEmb-RNN_Transfer"

$ws.Cells.Item(502, 1).Value = "outputs/2024-05-06/06-47-42"
$ws.Cells.Item(502, 2).Value = $true
$ws.Cells.Item(502, 3).Value = "mimiciii"
$ws.Cells.Item(502, 4).Value = "mortality"
$ws.Cells.Item(502, 5).Value = "DSVA_DPE"
$ws.Cells.Item(502, 6).Value = "descemb_rnn"
$ws.Cells.Item(502, 7).Value = "ehr_model"
$ws.Cells.Item(502, 8).Value = $false
$ws.Cells.Item(502, 9).Value = $false
$ws.Cells.Item(502, 10).Value = $false
$ws.Cells.Item(502, 11).Value = $true
$ws.Cells.Item(502, 13).Value = 92
$ws.Cells.Item(502, 14).Value = 1.053
$ws.Cells.Item(502, 15).Value = 0.507
$ws.Cells.Item(502, 16).Value = 0.09
$ws.Cells.Item(502, 17).Value = "DescEmb-RNN_Transfer"

$ws.Cells.Item(503, 1).Value = "outputs/2024-05-06/07-20-50"
$ws.Cells.Item(503, 2).Value = $true
$ws.Cells.Item(503, 3).Value = "eicu"
$ws.Cells.Item(503, 4).Value = "diagnosis"
$ws.Cells.Item(503, 5).Value = "DSVA_DPE"
$ws.Cells.Item(503, 6).Value = "descemb_rnn"
$ws.Cells.Item(503, 7).Value = "ehr_model"
$ws.Cells.Item(503, 8).Value = $false
$ws.Cells.Item(503, 9).Value = $false
$ws.Cells.Item(503, 10).Value = $false
$ws.Cells.Item(503, 11).Value = $true
$ws.Cells.Item(503, 13).Value = 218
$ws.Cells.Item(503, 14).Value = 0.634
$ws.Cells.Item(503, 15).Value = 0.794
$ws.Cells.Item(503, 16).Value = 0.454
$ws.Cells.Item(503, 17).Value = "DescEmb-RNN_Transfer"

$ws.Cells.Item(504, 1).Value = "outputs/2024-05-06/08-45-45"
$ws.Cells.Item(504, 2).Value = $true
$ws.Cells.Item(504, 3).Value = "eicu"
$ws.Cells.Item(504, 4).Value = "los_3day"
$ws.Cells.Item(504, 5).Value = "DSVA_DPE"
$ws.Cells.Item(504, 6).Value = "descemb_rnn"
$ws.Cells.Item(504, 7).Value = "ehr_model"
$ws.Cells.Item(504, 8).Value = $false
$ws.Cells.Item(504, 9).Value = $false
$ws.Cells.Item(504, 10).Value = $false
$ws.Cells.Item(504, 11).Value = $true
$ws.Cells.Item(504, 13).Value = 71
$ws.Cells.Item(504, 14).Value = 1.538
$ws.Cells.Item(504, 15).Value = 0.66
$ws.Cells.Item(504, 16).Value = 0.446
$ws.Cells.Item(504, 17).Value = "DescEmb-RNN_Transfer"

$ws.Cells.Item(505, 1).Value = "outputs/2024-05-06/09-12-32"
$ws.Cells.Item(505, 2).Value = $true
$ws.Cells.Item(505, 3).Value = "eicu"
$ws.Cells.Item(505, 4).Value = "los_7day"
$ws.Cells.Item(505, 5).Value = "DSVA_DPE"
$ws.Cells.Item(505, 6).Value = "descemb_rnn"
$ws.Cells.Item(505, 7).Value = "ehr_model"
$ws.Cells.Item(505, 8).Value = $false
$ws.Cells.Item(505, 9).Value = $false
$ws.Cells.Item(505, 10).Value = $false
$ws.Cells.Item(505, 11).Value = $true
$ws.Cells.Item(505, 13).Value = 88
$ws.Cells.Item(505, 14).Value = 1.44
$ws.Cells.Item(505, 15).Value = 0.638
$ws.Cells.Item(505, 16).Value = 0.153
$ws.Cells.Item(505, 17).Value = "DescEmb-RNN_Transfer"

$ws.Cells.Item(506, 1).Value = "outputs/2024-05-06/09-45-22"
$ws.Cells.Item(506, 2).Value = $true
$ws.Cells.Item(506, 3).Value = "eicu"
$ws.Cells.Item(506, 4).Value = "readmission"
$ws.Cells.Item(506, 5).Value = "DSVA_DPE"
$ws.Cells.Item(506, 6).Value = "descemb_rnn"
$ws.Cells.Item(506, 7).Value = "ehr_model"
$ws.Cells.Item(506, 8).Value = $false
$ws.Cells.Item(506, 9).Value = $false
$ws.Cells.Item(506, 10).Value = $false
$ws.Cells.Item(506, 11).Value = $true
$ws.Cells.Item(506, 13).Value = 47
$ws.Cells.Item(506, 14).Value = 1.026
$ws.Cells.Item(506, 15).Value = 0.497
$ws.Cells.Item(506, 16).Value = 0.105
$ws.Cells.Item(506, 17).Value = "DescEmb-RNN_Transfer"

$ws.Cells.Item(507, 1).Value = "outputs/2024-05-06/10-02-55"
$ws.Cells.Item(507, 2).Value = $true
$ws.Cells.Item(507, 3).Value = "eicu"
$ws.Cells.Item(507, 4).Value = "mortality"
$ws.Cells.Item(507, 5).Value = "DSVA_DPE"
$ws.Cells.Item(507, 6).Value = "descemb_rnn"
$ws.Cells.Item(507, 7).Value = "ehr_model"
$ws.Cells.Item(507, 8).Value = $false
$ws.Cells.Item(507, 9).Value = $false
$ws.Cells.Item(507, 10).Value = $false
$ws.Cells.Item(507, 11).Value = $true
$ws.Cells.Item(507, 13).Value = 141
$ws.Cells.Item(507, 14).Value = 1.248
$ws.Cells.Item(507, 15).Value = 0.666
$ws.Cells.Item(507, 16).Value = 0.144
$ws.Cells.Item(507, 17).Value = "DescEmb-RNN_Transfer"


# Column L holds the "patience" value but is stored as text in this sheet (e.g. "45"),
# matching the existing rows. Format the range as Text before writing so the numeric-
# looking string is not auto-coerced to a number, then drop back to the Normal style so
# no stray cell formatting is introduced.
$lRange = $ws.Range("L499:L507")
$lRange.NumberFormat = "@"
$lRange.Value = "45"
$lRange.Style = "Normal"
